$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = 10874
$ws.Range("C6:C7").Value = 10447
$ws.Range("C8:C12").Value = 9697
$ws.Range("C13:C25").Value = 9242
$ws.Range("C26").Value = 9131
$ws.Range("C27:C29").Value = 7793
$ws.Range("C30:C36").Value = 7345
$ws.Range("C37:C41").Value = 7312
$ws.Range("C42:C252").Value = 7310
